# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (fund-level breakdown) right after "2021-Q4"
#   and before "总计".
# - Re-create "总计" (so it gets bumped from sheetId 4 -> 5, matching the
#   target) with its original rows plus a new 2022-Q1 summary row on top.

$wb = $excel.ActiveWorkbook

# Reference sheet that already carries the header/index-column style (s="2")
# we want to replicate on the new sheets (border + bold + centered).
$styleRef = $wb.Worksheets.Item("2021-Q4")

# Capture the current "总计" data before we remove/recreate the sheet so we
# can rebuild it with the new row inserted on top.
$oldTotal = $wb.Worksheets.Item("总计")
$totalRows = @()
$r = 2
while ($oldTotal.Cells.Item($r, 2).Value() -ne $null) {
    $b = $oldTotal.Cells.Item($r, 2).Value()
    $c = $oldTotal.Cells.Item($r, 3).Value()
    $d = $oldTotal.Cells.Item($r, 4).Value()
    $totalRows += , @($b, $c, $d)
    $r = $r + 1
}

# Drop the old "总计" sheet -- its sheetId (4) will be reused by the sheet we
# insert next, which is what lets the recreated "总计" land on sheetId 5.
$oldTotal.Delete()

# ---------------------------------------------------------------------
# New sheet: 2022-Q1 (fund-level breakdown), inserted after "2021-Q4"
# ---------------------------------------------------------------------
$afterQ4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterQ4)
$q1.Name = "2022-Q1"

# Match the page margins used throughout the rest of this workbook
# (0.75in / 0.75in / 1in / 1in / 0.5in / 0.5in -> 54/54/72/72/36/36 points).
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Header row formatting (bold/border/center, style "2")
$styleRef.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1data = @(
    @("0", "159996", "国泰中证全指家用电器ETF", "17.91", "98.66", "2.24", "0.4012", 8),
    @("1", "005063", "广发中证全指家用电器指数A", "9.91", "94.24", "2.14", "0.2121", 9),
    @("2", "005064", "广发中证全指家用电器指数C", "5.40", "94.24", "2.14", "0.1156", 9),
    @("3", "005041", "人保研究精选混合A", "1.33", "81.87", "2.21", "0.0294", 4),
    @("4", "006573", "人保行业轮动混合A", "0.97", "81.10", "2.58", "0.0250", 2),
    @("5", "006574", "人保行业轮动混合C", "0.18", "81.10", "2.58", "0.0046", 2),
    @("6", "005042", "人保研究精选混合C", "0.03", "81.87", "2.21", "0.0007", 4)
)

$row = 2
foreach ($rec in $q1data) {
    # Index column (A) carries the same bold/border/center style as the header.
    $styleRef.Range("A2").Copy()
    $q1.Cells.Item($row, 1).PasteSpecial(-4122)
    $q1.Cells.Item($row, 1).Value = [int]$rec[0]

    # Fund code (keep leading zeros -- e.g. "005063" -- by forcing text).
    $q1.Cells.Item($row, 2).NumberFormat = "@"
    $q1.Cells.Item($row, 2).Value = $rec[1]

    $q1.Cells.Item($row, 3).Value = $rec[2]

    # D-G are stored as text in the source data (keep them as text so they
    # don't get auto-coerced to plain numbers).
    $q1.Cells.Item($row, 4).NumberFormat = "@"
    $q1.Cells.Item($row, 4).Value = $rec[3]

    $q1.Cells.Item($row, 5).NumberFormat = "@"
    $q1.Cells.Item($row, 5).Value = $rec[4]

    $q1.Cells.Item($row, 6).NumberFormat = "@"
    $q1.Cells.Item($row, 6).Value = $rec[5]

    $q1.Cells.Item($row, 7).NumberFormat = "@"
    $q1.Cells.Item($row, 7).Value = $rec[6]

    $q1.Cells.Item($row, 8).Value = [int]$rec[7]

    $row = $row + 1
}

# ---------------------------------------------------------------------
# Recreate "总计" right after the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$styleRef.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# New row for 2022-Q1, followed by the rows that existed before the edit.
$allTotalRows = New-Object System.Collections.ArrayList
[void]$allTotalRows.Add(@("2022-Q1", 7, 0.79))
foreach ($old in $totalRows) {
    [void]$allTotalRows.Add($old)
}

$row = 2
$idx = 0
foreach ($rec in $allTotalRows) {
    $styleRef.Range("A2").Copy()
    $total.Cells.Item($row, 1).PasteSpecial(-4122)
    $total.Cells.Item($row, 1).Value = $idx

    $total.Cells.Item($row, 2).Value = $rec[0]
    $total.Cells.Item($row, 3).Value = $rec[1]
    $total.Cells.Item($row, 4).Value = $rec[2]

    $row = $row + 1
    $idx = $idx + 1
}
